$d = $word.ActiveDocument

# The document's headers/footers each carry a single inline logo picture.
# Per the commit, the BTec logo picture is renamed image1.jpg -> image2.jpg,
# and the two Pearson logo pictures (one in the default footer, one in the
# first-page footer) are each renamed image2.png -> image1.png.
#
# The InlineShape.Name *getter* in this host doesn't reflect the picture's
# existing docPr/@name, so shapes are identified by their (stable)
# AlternativeText/description instead, and renamed via the Name *setter*
# (which does write through to the underlying drawing XML).

foreach ($sec in $d.Sections) {
    for ($i = 1; $i -le 3; $i++) {
        $hdr = $sec.Headers.Item($i)
        if ($hdr.Exists) {
            foreach ($shp in $hdr.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                }
            }
        }

        $ftr = $sec.Footers.Item($i)
        if ($ftr.Exists) {
            foreach ($shp in $ftr.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }
}
